$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values (column D) read as plain numbers to Excel's
# automatic type detection (e.g. "226.09"), whereas in the source workbook
# every Price cell is stored as text (note values like "34.303.01" that
# aren't valid numbers at all). Force those specific cells to Text format
# first so their new values are written back out as text too.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "34.303.01"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.790.74"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "226.09"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "32.63"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "0.295"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "0.0689"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "1.790.51"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "0.632"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "34.296.45"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "4.27"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "68.35"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "243.61"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "11.29"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "165.73"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").Value = "7.29"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").Value = "16.47"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "3.96"
$ws.Range("E30").Value = "  +6.08%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").Value = "1.399.48"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("D37").Value = "0.666"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "85.08"
$ws.Range("E40").Value = "  +3.93%  "
$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  +4.03%  "
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "0.935"
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("D44").Value = "13.70"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").Value = "0.0524"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").Value = "1.949.08"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "104.67"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  -1.28%  "
